$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# New method / note entries appended below the existing list (rows 11-15, column A)
$ws.Range("A11").Value = "MessageManager.deleteMessage(Message message);"
$ws.Range("A12").Value = "get membre by nickname(nickname)"
$ws.Range("A13").Value = "get membre by nickname et password(nickname+password)"
$ws.Range("A14").Value = "getClinsdoeilRecus(int TOId){"
$ws.Range("A15").Value = "ClinsdoeilManager.deleteClindoeil(Clinsdoeil clin){"

# Rows 11-14 use the existing red-font style (same as the other column A entries, s="1")
$ws.Range("A11:A14").Font.Color = 255

# Update the selected cell to reflect the new last populated cell
$ws.Range("C15").Select()
